$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$ws.Range("A2").Copy()
$scratch.PasteSpecial(-4122)
$excel.CutCopyMode = 0
$scratch.Interior.Pattern = -4142
Write-Host "fg:" $scratch.Interior.Color "bg:" $scratch.Interior.PatternColor
